# "最终消费" (final consumption) sheet update:
#   - Drop the four oldest year rows (2000年, 2002年, 2005年, 2007年).
#   - That shifts the existing 2010年/2012年/2015年/2017年 rows up from
#     rows 6-9 to rows 2-5.
#   - Append a brand new 2020年 row as the new row 6.
#   - The used range shrinks from A1:S9 to A1:S6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting rows 2-5 (2000年..2007年) shifts rows 6-9 (2010年..2017年) up to
# rows 2-5, preserving all of their original values/styles/blank cells.
$ws.Range("A2:A5").EntireRow.Delete() | Out-Null

# Row 6 is now empty; fill it in with the 2020年 figures.
$year2020 = @{
    1  = "2020年"
    2  = 4554576.12250768
    3  = 2324055296.5576
    4  = 350440083.406576
    9  = 5576336944.0057
    10 = 276123836.16857
    11 = 126886993.017075
    12 = 116095264.790195
    13 = 131207526.795881
    15 = 1222279.57958158
    16 = 4880487.07462051
    18 = 5486655.16723377
    19 = 608515165.754948
}
# Columns E, F, G, H, N, Q (5,6,7,8,14,17) have no reported value for 2020年
# and are intentionally left blank, matching the source data.

foreach ($col in $year2020.Keys) {
    $ws.Cells.Item(6, $col).Value = $year2020[$col]
}

# Give the new 2020年 label (column A) the same look as the other year
# labels in column A: bold, centered/top-aligned, thin box border.
$yearCell = $ws.Cells.Item(6, 1)
$yearCell.Font.Bold = $true
$yearCell.HorizontalAlignment = -4108   # xlCenter
$yearCell.VerticalAlignment = -4160     # xlTop
$yearCell.Borders.LineStyle = 1         # xlContinuous (thin box border)
